$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.879.20"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.229.16"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.95"
$ws.Range("E5").Value = "  +6.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.02"
$ws.Range("E7").Value = "  +3.31%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +5.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.22"
$ws.Range("E10").Value = "  +14.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.11"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.15"
$ws.Range("E13").Value = "  +5.78%  "
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "2.561.02"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.99"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.865"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "2.230.63"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "41.879.86"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "0.0₃0967"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.89"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.42"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("E24").Value = "  +5.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.11"
$ws.Range("E25").Value = "  +12.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +4.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.69"
$ws.Range("E28").Value = "  +6.25%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.34"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.76"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.125"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.56"
$ws.Range("E34").Value = "  +2.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0732"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.71"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.21"
$ws.Range("E37").Value = "  +20.99%  "
$ws.Range("E38").Value = "  +8.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0302"
$ws.Range("E39").Value = "  +12.13%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.28"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.98"
$ws.Range("E43").Value = "  +18.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.95"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.202"
$ws.Range("E45").Value = "  +6.25%  "
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.67"
$ws.Range("E48").Value = "  +5.30%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  +6.72%  "
$ws.Range("E51").Value = "  +0.99%  "
